$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (92-100) to append below the existing data (which ends at row 91)
$newRows = @(
    @(45725.733101851853, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.737037037034, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.739548611113, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.742303240739, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.743796296294, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.769189814811, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.784131944441, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.784699074073, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686),
    @(45725.786273148151, 10, 6, 240, 423, 399, 476, 3432, 476, 2026, 208, 417, 30, 3667, 4686)
)

$startRow = 92
$lastTemplateRow = 91

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Copy formatting from the last existing row so number formats / styles match
    $srcRange = $ws.Range($ws.Cells.Item($lastTemplateRow, 1), $ws.Cells.Item($lastTemplateRow, 15))
    $dstRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 15))
    $srcRange.Copy($dstRange)

    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
